$d = $word.ActiveDocument

# 1. Title heading and the later bold duplicate of the same text
$d.Content.Find.Execute(
    "Play Jellyfish Flow for Free - A High-Volatility Slot Game", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Jellyfish Flow for Free", 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute(
    "Stunning oceanic graphics and cheerful musical theme", $true, $false, $false, $false, $false,
    $true, 1, $false, "Expanding reels and ways to win", 2)

$d.Content.Find.Execute(
    "Unique game mechanics with reel expansion and wilds", $true, $false, $false, $false, $false,
    $true, 1, $false, "Two different volatility options", 2)

$d.Content.Find.Execute(
    "Option to choose between regular and ultra-high volatility modes", $true, $false, $false, $false, $false,
    $true, 1, $false, "Beautiful graphics and theme", 2)

$d.Content.Find.Execute(
    "Massive potential payouts of up to 32,416x or 168,004x", $true, $false, $false, $false, $false,
    $true, 1, $false, "High potential payouts", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute(
    "Not recommended for casual players or beginners", $true, $false, $false, $false, $false,
    $true, 1, $false, "Varied RTP depending on online casino", 2)

$d.Content.Find.Execute(
    "Actual RTP is subject to variation depending on the online casino", $true, $false, $false, $false, $false,
    $true, 1, $false, "Not suitable for casual players or beginners", 2)

# 4. Italic meta description near the end
$d.Content.Find.Execute(
    "Read a review of Jellyfish Flow slot game, play it for free, and experience its unique reel expansion and wilds mechanics. Choose your volatility mode and aim for massive payouts.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Jellyfish Flow and play this high-risk, high-reward slot game for free.", 2)
